# "Updated credentials of VM in the poster"
#
# The poster has a "Virtual Machine Settings" section with a "Remote Access"
# box (shape "CustomShape 23") showing a Partner ID / Password pair for
# remote-desktop access to the project VM. This edit refreshes the Partner
# ID value and bumps its displayed size to match the Password value's size,
# and it also turns the "Ontmalizer" tool link (shape "CustomShape 33") into
# a normal-looking hyperlink (blue + underlined) instead of plain black text.
# Finally it clears the stock "Click to edit the title text format" prompt
# left over on the slide master's (unused) title placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Remote Access box: update the Partner ID and match its font size to
#    the rest of the credential values (28pt -> 32pt).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        if ($full -ne $null -and $full.Length -gt 0) {
            $needle = "905 399 395"
            $idx = $full.IndexOf($needle)
            if ($idx -ne $null -and $idx -ge 0) {
                $sub = $tr.Characters($idx + 1, $needle.Length)
                $sub.Font.Size = 32
                $sub.Text = "873 887 809"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) "Ontmalizer" reference link: style it like a hyperlink (blue,
#    underlined) instead of plain black text.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        if ($full -ne $null -and $full.Length -gt 0) {
            $needle = "http://github.com/srdc/ontmalizer"
            $idx = $full.IndexOf($needle)
            if ($idx -ne $null -and $idx -ge 0) {
                $sub = $tr.Characters($idx + 1, $needle.Length)
                $sub.Font.Underline = 1
                $sub.Font.Color.RGB = 16711680
            }
        }
    }
}
